$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,"BTC","Bitcoin",66767,1317540462031,58370774815,-2.30488),
    @(2,"ETH","Ethereum",3570.76,429377428195,26843432307,-3.69067),
    @(3,"USDT","Tether",0.996949,103211906265,73693896831,-0.56706),
    @(4,"BNB","BNB",570.55,88100504532,4885730188,-5.72398),
    @(5,"SOL","Solana",191.34,85203844660,12598405929,0.129),
    @(6,"STETH","Lido Staked Ether",3567.87,34946620795,76002939,-3.50899),
    @(7,"XRP","XRP",0.613313,33627532063,2418528876,-2.48676),
    @(8,"USDC","USDC",0.998122,30780165989,11918839584,-0.15399),
    @(9,"ADA","Cardano",0.676299,23921770520,1068617556,-5.60821),
    @(10,"DOGE","Dogecoin",0.148312,21340672957,3020361808,-5.39815),
    @(11,"AVAX","Avalanche",55.83,21173856921,1990840741,-4.93622),
    @(12,"SHIB","Shiba Inu",0.0000268,15854269859,2344238830,-4.63066),
    @(13,"DOT","Polkadot",9.83,13214172041,583278730,-4.83464),
    @(14,"TON","Toncoin",3.62,12591205708,105687801,-2.64036),
    @(15,"TRX","TRON",0.124773,10968818908,515129095,-1.81144),
    @(16,"LINK","Chainlink",18.24,10771861912,794227424,-5.41016),
    @(17,"WBTC","Wrapped Bitcoin",66936,10399368080,418936379,-1.8624),
    @(18,"MATIC","Polygon",1.059,9863606995,892614640,-7.02731),
    @(19,"UNI","Uniswap",12.12,9147147043,267528590,-4.67625),
    @(20,"NEAR","NEAR Protocol",7.76,8051637249,1068358670,1.00923),
    @(21,"BCH","Bitcoin Cash",399.14,7877840183,593591227,-1.49163),
    @(22,"LTC","Litecoin",85.81999999999999,6388271232,720251547,-4.07834),
    @(23,"ICP","Internet Computer",12.42,5737850176,211363273,-3.12471),
    @(24,"LEO","LEO Token",6.03,5600043458,3610533,0.06342),
    @(25,"APT","Aptos",14.08,5594377088,534345126,-5.66227),
    @(26,"ATOM","Cosmos Hub",12.1,4740189264,368432907,-3.69979),
    @(27,"DAI","Dai",0.998034,4683938253,768256322,-0.11995),
    @(28,"FIL","Filecoin",8.890000000000001,4679781583,527310065,-6.51051),
    @(29,"ETC","Ethereum Classic",31.07,4544696452,312339699,-3.89611),
    @(30,"RNDR","Render",11.57,4414049385,665371836,-1.45818),
    @(31,"TAO","Bittensor",638.78,4118602027,31884705,2.74366),
    @(32,"IMX","Immutable",2.91,4052562347,106967040,-3.60499),
    @(33,"STX","Stacks",2.72,3938353274,128412337,-0.45104),
    @(34,"HBAR","Hedera",0.113801,3838591906,99590958,-5.03746),
    @(35,"OKB","OKB",63.85,3832538080,27764541,-5.51109),
    @(36,"GRT","The Graph",0.399825,3770222674,269482757,-1.59098),
    @(37,"INJ","Injective",42.03,3720439441,366997991,-9.737579999999999),
    @(38,"XLM","Stellar",0.129542,3716780240,172822366,-6.21046),
    @(39,"CRO","Cronos",0.138837,3698253250,41159725,-4.25631),
    @(40,"OP","Optimism",3.57,3601601450,370424275,-4.43856),
    @(41,"PEPE","Pepe",0.00000761,3214387150,1696661353,-5.37792),
    @(42,"KAS","Kaspa",0.133093,3071534914,68019190,-2.40645),
    @(43,"VET","VeChain",0.04129212,3010039080,131567142,-5.54854),
    @(44,"FDUSD","First Digital USD",0.992972,2998597438,14657445902,-0.77877),
    @(45,"MKR","Maker",3221.46,2974708156,353552149,12.16911),
    @(46,"THETA","Theta Network",2.95,2953977853,98078283,-0.4347),
    @(47,"FET","Fetch.ai",2.7,2824030564,688189129,3.89427),
    @(48,"MNT","Mantle",0.84972,2750073339,381630426,-1.72282),
    @(49,"XMR","Monero",142.51,2592686166,41985615,-1.85836),
    @(50,"RUNE","THORChain",8.58,2579222805,705667367,-5.75554)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    $ws.Range("A$row").Value = $rowValues[0]
    $ws.Range("B$row").Value = $rowValues[1]
    $ws.Range("C$row").Value = $rowValues[2]
    $ws.Range("D$row").Value = $rowValues[3]
    $ws.Range("E$row").Value = $rowValues[4]
    $ws.Range("F$row").Value = $rowValues[5]
    $ws.Range("G$row").Value = $rowValues[6]
}
